# Refresh the cryptocurrency table on Sheet1: updates the "Price" (column D)
# and "Volume(1h)" (column E) figures with newly scraped values, matching
# the scheduled "Updated cryptos list ... with GitHub Actions" run.
#
# Both columns are stored as plain text (prices sometimes use "." as a
# thousands separator, e.g. "63.289.55", and percentages keep their padding
# spaces, e.g. "  +3.07%  "), so plain-decimal values such as "593.87" are
# written with a forced text format -- otherwise Excel's normal
# type-inference would silently turn them into numbers, which would change
# the cell type from the original text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '63.289.55' },
    @{ Cell = 'D3'; Value = '3.240.80' },
    @{ Cell = 'E3'; Value = '  +3.07%  ' },
    @{ Cell = 'E4'; Value = '  -0.04%  ' },
    @{ Cell = 'D5'; Value = '593.87' },
    @{ Cell = 'E5'; Value = '  -1.27%  ' },
    @{ Cell = 'D6'; Value = '140.49' },
    @{ Cell = 'E6'; Value = '  -1.01%  ' },
    @{ Cell = 'E7'; Value = '  -0.02%  ' },
    @{ Cell = 'D8'; Value = '3.236.61' },
    @{ Cell = 'E9'; Value = '  -1.72%  ' },
    @{ Cell = 'E10'; Value = '  -0.83%  ' },
    @{ Cell = 'E11'; Value = '  -0.44%  ' },
    @{ Cell = 'E12'; Value = '  -0.29%  ' },
    @{ Cell = 'E13'; Value = '  -2.63%  ' },
    @{ Cell = 'D14'; Value = '34.31' },
    @{ Cell = 'E14'; Value = '  -1.53%  ' },
    @{ Cell = 'D15'; Value = '3.770.58' },
    @{ Cell = 'E15'; Value = '  +2.90%  ' },
    @{ Cell = 'E16'; Value = '  -0.18%  ' },
    @{ Cell = 'D17'; Value = '3.237.67' },
    @{ Cell = 'E17'; Value = '  +3.20%  ' },
    @{ Cell = 'D18'; Value = '63.314.11' },
    @{ Cell = 'E18'; Value = '  -1.00%  ' },
    @{ Cell = 'E19'; Value = '  -1.04%  ' },
    @{ Cell = 'D20'; Value = '474.01' },
    @{ Cell = 'E20'; Value = '  -2.55%  ' },
    @{ Cell = 'D21'; Value = '14.15' },
    @{ Cell = 'E21'; Value = '  -3.53%  ' },
    @{ Cell = 'D22'; Value = '0.731' },
    @{ Cell = 'E22'; Value = '  +2.84%  ' },
    @{ Cell = 'E23'; Value = '  +2.87%  ' },
    @{ Cell = 'D24'; Value = '83.88' },
    @{ Cell = 'E24'; Value = '  -5.09%  ' },
    @{ Cell = 'D25'; Value = '13.18' },
    @{ Cell = 'E25'; Value = '  -0.26%  ' },
    @{ Cell = 'E27'; Value = '  -0.99%  ' },
    @{ Cell = 'D28'; Value = '7.36' },
    @{ Cell = 'D29'; Value = '8.10' },
    @{ Cell = 'E29'; Value = '  -0.94%  ' },
    @{ Cell = 'D30'; Value = '2.12' },
    @{ Cell = 'E30'; Value = '  +3.01%  ' },
    @{ Cell = 'D31'; Value = '27.50' },
    @{ Cell = 'E31'; Value = '  +0.29%  ' },
    @{ Cell = 'E32'; Value = '  -0.08%  ' },
    @{ Cell = 'E33'; Value = '  -3.85%  ' },
    @{ Cell = 'E34'; Value = '  -4.42%  ' },
    @{ Cell = 'E35'; Value = '  -1.22%  ' },
    @{ Cell = 'D36'; Value = '5.92' },
    @{ Cell = 'E36'; Value = '  -1.95%  ' },
    @{ Cell = 'D37'; Value = '52.67' },
    @{ Cell = 'E37'; Value = '  -0.02%  ' },
    @{ Cell = 'E38'; Value = '  -4.18%  ' },
    @{ Cell = 'D39'; Value = '0.0393' },
    @{ Cell = 'E39'; Value = '  -1.00%  ' },
    @{ Cell = 'D40'; Value = '423.07' },
    @{ Cell = 'E40'; Value = '  -1.81%  ' },
    @{ Cell = 'D41'; Value = '8.38' },
    @{ Cell = 'E41'; Value = '  +0.23%  ' },
    @{ Cell = 'D42'; Value = '2.972.93' },
    @{ Cell = 'E42'; Value = '  +2.14%  ' },
    @{ Cell = 'E43'; Value = '  -5.83%  ' },
    @{ Cell = 'E44'; Value = '  -8.24%  ' },
    @{ Cell = 'D45'; Value = '0.267' },
    @{ Cell = 'E45'; Value = '  +2.80%  ' },
    @{ Cell = 'E46'; Value = '  -0.69%  ' },
    @{ Cell = 'E47'; Value = '  +0.06%  ' },
    @{ Cell = 'E48'; Value = '  +0.66%  ' },
    @{ Cell = 'D49'; Value = '2.32' },
    @{ Cell = 'E49'; Value = '  -2.80%  ' },
    @{ Cell = 'E50'; Value = '  -0.35%  ' },
    @{ Cell = 'D51'; Value = '121.44' },
    @{ Cell = 'E51'; Value = '  +0.36%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)

    if ($u.Value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Plain decimal number text (e.g. "593.87") -- force text storage
        # (equivalent to an apostrophe-prefixed entry in the Excel UI) so it
        # keeps its original string cell type instead of becoming a number.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
